$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# out_vars: append row 11 (new date 2020-06-10 / serial 43992) by cloning the
# formatting of row 10 and writing the new figures.
# ---------------------------------------------------------------------------
$wsOut = $wb.Worksheets.Item("out_vars")

$wsOut.Range("A10:J10").Copy()
$wsOut.Range("A11:J11").PasteSpecial(-4122)

$wsOut.Range("A11").Value = 43992
$wsOut.Range("B11").Value = 129184
$wsOut.Range("C11").Value = 186570
$wsOut.Range("D11").Value = 53608
$wsOut.Range("E11").Value = 15357
$wsOut.Range("F11").Value = 33.110137478325491
$wsOut.Range("G11").Value = 42773
$wsOut.Range("H11").Value = 3970
$wsOut.Range("I11").Value = 4126
$wsOut.Range("J11").Value = 369362

$wsOut.Range("I13").Select()

# ---------------------------------------------------------------------------
# dates_dx: row 11 placeholders already exist (formatted, blank) - fill them.
# ---------------------------------------------------------------------------
$wsDx = $wb.Worksheets.Item("dates_dx")

$wsDx.Range("A11").Value = 43992
$wsDx.Range("B11").Value = 0
$wsDx.Range("C11").Value = 1
$wsDx.Range("D11").Value = 1
$wsDx.Range("E11").Value = 1
$wsDx.Range("F11").Value = 0
$wsDx.Range("G11").Value = 0
$wsDx.Range("H11").Value = 0
$wsDx.Range("I11").Value = 4

$wsDx.Range("J11").Select()

# ---------------------------------------------------------------------------
# dates_sx: no row 11 yet - clone the date style from A10 and add the row.
# ---------------------------------------------------------------------------
$wsSx = $wb.Worksheets.Item("dates_sx")

$wsSx.Range("A10").Copy()
$wsSx.Range("A11").PasteSpecial(-4122)

$wsSx.Range("A11").Value = 43992
$wsSx.Range("B11").Value = 0
$wsSx.Range("C11").Value = 1
$wsSx.Range("D11").Value = 0
$wsSx.Range("E11").Value = 1
$wsSx.Range("F11").Value = 1
$wsSx.Range("G11").Value = 1
$wsSx.Range("H11").Value = 0
$wsSx.Range("I11").Value = 1
$wsSx.Range("J11").Value = 1
$wsSx.Range("K11").Value = 0
$wsSx.Range("L11").Value = 0

$wsSx.Range("L11").Select()

# ---------------------------------------------------------------------------
# dates_deaths: row 11 existed as a blank placeholder (A11 only) - populate it
# with the same look as row 10.
# ---------------------------------------------------------------------------
$wsDeaths = $wb.Worksheets.Item("dates_deaths")

$wsDeaths.Range("A10").Copy()
$wsDeaths.Range("A11").PasteSpecial(-4122)

$wsDeaths.Range("A11").Value = 43992
$wsDeaths.Range("B11").Value = 1
$wsDeaths.Range("C11").Value = 0
$wsDeaths.Range("D11").Value = 2
$wsDeaths.Range("E11").Value = 1
$wsDeaths.Range("F11").Value = 1
$wsDeaths.Range("G11").Value = 2
$wsDeaths.Range("H11").Value = 2

$wsDeaths.Range("I11").Select()

# ---------------------------------------------------------------------------
# control_obs: new date column K (2020-06-10) across the existing metric
# rows, plus extending the running-total SUM formulas into J20:K20.
# ---------------------------------------------------------------------------
$wsControl = $wb.Worksheets.Item("control_obs")

$wsControl.Range("K1").Value = 43992
$wsControl.Range("K2").Value = 3337
$wsControl.Range("K3").Value = 3152
$wsControl.Range("K4").Value = 3152
$wsControl.Range("K5").Value = 3152
$wsControl.Range("K6").Value = 3152
$wsControl.Range("K7").Value = 2449
$wsControl.Range("K8").Value = 4974
$wsControl.Range("K10").Value = 150
$wsControl.Range("K11").Value = 150
$wsControl.Range("K12").Value = 150
$wsControl.Range("K13").Value = 150
$wsControl.Range("K14").Value = 150
$wsControl.Range("K15").Value = 127
$wsControl.Range("K16").Value = 162
$wsControl.Range("K18").Value = 806

$wsControl.Range("J20").Formula = "=SUM(J2:J18)"
$wsControl.Range("K20").Formula = "=SUM(K2:K18)"

$wsControl.Range("K23").Select()

Write-Output "applied bitacora_historica_datos_abiertos update"
